$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.250.13"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.855.90"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.23"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4743"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2745"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06419"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").Value = "1.857.02"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07425"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.09"
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.979"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.35"
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6321"
$ws.Range("E15").Value = "  -4.24%  "
$ws.Range("D16").Value = "30.223.43"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.75"
$ws.Range("E18").Value = "  -4.07%  "
$ws.Range("E19").Value = "  -3.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "224.82"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "2.090.34"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.097"
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.028"
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.28"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.222"
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.76"
$ws.Range("E27").Value = "  -3.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.854"
$ws.Range("E28").Value = "  -5.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1022"
$ws.Range("E29").Value = "  +9.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.381"
$ws.Range("E30").Value = "  -5.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.223"
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.897"
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04882"
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.146"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7255"
$ws.Range("E35").Value = "  -2.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9997"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("E38").Value = "  +5.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.626"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8989"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.974"
$ws.Range("E41").Value = "  -4.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.80"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9940"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4092"
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.523"
$ws.Range("E45").Value = "  -6.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.027"
$ws.Range("E46").Value = "  -5.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.16"
$ws.Range("E47").Value = "  -5.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1203"
$ws.Range("E48").Value = "  -5.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.785"
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.400"
$ws.Range("E50").Value = "  -5.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05587"
$ws.Range("E51").Value = "  -0.83%  "
